$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added. Insert a new row at position 4,
# which pushes the existing rows 4-31 down to 5-32 (row 32 becomes the
# former row 31's data).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new record's values.
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Vega Modelo de Temuco"
$ws.Range("C4").Value = "La Araucanía"
$ws.Range("D4").Value = 44699
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 100112010
$ws.Range("G4").Value = "Achicoria"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 110
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 12000
$ws.Range("N4").Value = "$/caja 18 unidades"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 667
$ws.Range("Q4").Value = 18
$ws.Range("R4").Value = "Hortaliza"
